# Add a new "CEGEP" part-number column (column B) in front of the existing
# data on each of the three sheets (Cap, Res, Ind).

$wb = $excel.ActiveWorkbook

$sheetInfo = @(
    @{ Name = "Cap"; FilterRef = "A1:S1" },
    @{ Name = "Res"; FilterRef = "A1:R15" },
    @{ Name = "Ind"; FilterRef = "A1:R1" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Insert a new blank column before the current column B ("Value"/...).
    $ws.Columns.Item(2).Insert()

    # New header + label for the inserted column.
    $ws.Cells.Item(1, 2).Value = "CEGEP"

    # Re-apply the AutoFilter so it covers the new column.
    $ws.AutoFilterMode = $false
    $ws.Range($info.FilterRef).AutoFilter()
}

# The hidden _FilterDatabase defined names keep the autofilter's pre-insert
# range; update them to match the new extents.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Cap!_FilterDatabase") {
        $n.RefersTo = "=Cap!`$A`$1:`$S`$1"
    } elseif ($n.Name -eq "Res!_FilterDatabase") {
        $n.RefersTo = "=Res!`$A`$1:`$R`$15"
    } elseif ($n.Name -eq "Ind!_FilterDatabase") {
        $n.RefersTo = "=Ind!`$A`$1:`$R`$1"
    }
}

# Make "Cap" the active/selected sheet (it was the sheet being edited).
$wb.Worksheets.Item("Cap").Select()
